$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: I1 = "I0", J1 = "IF"
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the existing header formatting (bold font, border, centered/top alignment)
# from H1 onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-16 for the two new columns I and J
$data = @(
    @(5, 5),
    @(8, 8),
    @(7, 7),
    @(9, 9),
    @(7, 7),
    @(7, 8),
    @(8, 8),
    @(7, 7),
    @(9, 9),
    @(7, 7),
    @(6, 6),
    @(3, 6),
    @(1, 4),
    @(6, 7),
    @(1, 2)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
